$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 ("노브랜드"): confirmed offering price (D) went from unset "-" to 14000,
# and the offering amount (E) was revised from 10440 to 16800.
# Both columns store their numbers as shared-string text (not numeric cells) in
# this workbook, so force text entry with a leading apostrophe and then restore
# the default "Normal" style so no stray number-format style is left behind.

$ws.Range("D15").Value = "'14000"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "'16800"
$ws.Range("E15").Style = "Normal"
